$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.269999999999999
$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.452
$ws.Range("C13").Value = -12.694
$ws.Range("A18").Value = -21.766
$ws.Range("D20").Value = -8.222
$ws.Range("A25").Value = -21.534
